$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 values (Rakesh* -> Automation* naming, SF - Store -> parola)
$ws.Range("D2").Value = "AutomationModel"
$ws.Range("F2").Value = "AutomationAsset"
$ws.Range("H2").Value = "AutomationSupplier"
$ws.Range("M2").Value = "parola"

# Remove row 3 (Rakesh22 / us-9877 / ASAN entry) entirely
$ws.Rows.Item(3).Delete()

# Widen column D and add an explicit width for the newly meaningful column H
$ws.Columns.Item(4).ColumnWidth = 19
$ws.Columns.Item(8).ColumnWidth = 20.83

# Move the active selection to C5
$ws.Range("C5").Select()
